$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = [double]"0.008135666666666668"
$ws.Range("H2").Value = [double]"0.024407"
$ws.Range("I2").Value = [double]"0.0001175588769867851"
$ws.Range("J2").Value = [double]"0.0001175588769867851"
$ws.Range("M2").Value = [double]"0.1293233333333333"
$ws.Range("N2").Value = [double]"0.38797"
$ws.Range("O2").Value = [double]"0.02793463826022293"
$ws.Range("P2").Value = [double]"0.02793463826022293"
$ws.Range("Q2").Value = [double]"0.001052131532222222"
$ws.Range("R2").Value = [double]"0.009469183790000002"
$ws.Range("S2").Value = [double]"3.283964702903889E-06"
$ws.Range("T2").Value = [double]"3.283964702903888E-06"

$ws.Range("G3").Value = [double]"0.008135666666666668"
$ws.Range("H3").Value = [double]"0.024407"
$ws.Range("I3").Value = [double]"0.0001175588769867851"
$ws.Range("J3").Value = [double]"0.0001175588769867851"
$ws.Range("M3").Value = [double]"0.9720173333333334"
$ws.Range("N3").Value = [double]"2.916052"
$ws.Range("O3").Value = [double]"0.2099617438667928"
$ws.Range("P3").Value = [double]"0.2099617438667928"
$ws.Range("Q3").Value = [double]"0.007908009018222225"
$ws.Range("R3").Value = [double]"0.07117208116400001"
$ws.Range("S3").Value = [double]"2.468286681916718E-05"
$ws.Range("T3").Value = [double]"2.468286681916717E-05"

$ws.Range("G4").Value = [double]"0.008135666666666668"
$ws.Range("H4").Value = [double]"0.024407"
$ws.Range("I4").Value = [double]"0.0001175588769867851"
$ws.Range("J4").Value = [double]"0.0001175588769867851"
$ws.Range("M4").Value = [double]"1.000161333333333"
$ws.Range("N4").Value = [double]"3.000484"
$ws.Range("O4").Value = [double]"0.2160410215882329"
$ws.Range("P4").Value = [double]"0.2160410215882329"
$ws.Range("Q4").Value = [double]"0.00813697922088889"
$ws.Range("R4").Value = [double]"0.073232812988"
$ws.Range("S4").Value = [double]"2.539753988099047E-05"
$ws.Range("T4").Value = [double]"2.539753988099046E-05"

$ws.Range("G5").Value = [double]"0.008135666666666668"
$ws.Range("H5").Value = [double]"0.024407"
$ws.Range("I5").Value = [double]"0.0001175588769867851"
$ws.Range("J5").Value = [double]"0.0001175588769867851"
$ws.Range("M5").Value = [double]"2.527995333333333"
$ws.Range("N5").Value = [double]"7.583985999999999"
$ws.Range("O5").Value = [double]"0.5460625962847514"
$ws.Range("P5").Value = [double]"0.5460625962847514"
$ws.Range("Q5").Value = [double]"0.02056692736688889"
$ws.Range("R5").Value = [double]"0.185102346302"
$ws.Range("S5").Value = [double]"6.419450558372361E-05"
$ws.Range("T5").Value = [double]"6.41945055837236E-05"

$ws.Range("G6").Value = [double]"0.01352566666666667"
$ws.Range("H6").Value = [double]"0.040577"
$ws.Range("I6").Value = [double]"0.0001954433790098242"
$ws.Range("J6").Value = [double]"0.0001954433790098242"
$ws.Range("M6").Value = [double]"0.1293233333333333"
$ws.Range("N6").Value = [double]"0.38797"
$ws.Range("O6").Value = [double]"0.02793463826022293"
$ws.Range("P6").Value = [double]"0.02793463826022293"
$ws.Range("Q6").Value = [double]"0.001749184298888889"
$ws.Range("R6").Value = [double]"0.01574265869"
$ws.Range("S6").Value = [double]"5.459640092995087E-06"
$ws.Range("T6").Value = [double]"5.459640092995086E-06"

$ws.Range("G7").Value = [double]"0.01352566666666667"
$ws.Range("H7").Value = [double]"0.040577"
$ws.Range("I7").Value = [double]"0.0001954433790098242"
$ws.Range("J7").Value = [double]"0.0001954433790098242"
$ws.Range("M7").Value = [double]"0.9720173333333334"
$ws.Range("N7").Value = [double]"2.916052"
$ws.Range("O7").Value = [double]"0.2099617438667928"
$ws.Range("P7").Value = [double]"0.2099617438667928"
$ws.Range("Q7").Value = [double]"0.01314718244488889"
$ws.Range("R7").Value = [double]"0.118324642004"
$ws.Range("S7").Value = [double]"4.103563268412121E-05"
$ws.Range("T7").Value = [double]"4.103563268412121E-05"

$ws.Range("G8").Value = [double]"0.01352566666666667"
$ws.Range("H8").Value = [double]"0.040577"
$ws.Range("I8").Value = [double]"0.0001954433790098242"
$ws.Range("J8").Value = [double]"0.0001954433790098242"
$ws.Range("M8").Value = [double]"1.000161333333333"
$ws.Range("N8").Value = [double]"3.000484"
$ws.Range("O8").Value = [double]"0.2160410215882329"
$ws.Range("P8").Value = [double]"0.2160410215882329"
$ws.Range("Q8").Value = [double]"0.01352784880755555"
$ws.Range("R8").Value = [double]"0.121750639268"
$ws.Range("S8").Value = [double]"4.222378726393862E-05"
$ws.Range("T8").Value = [double]"4.222378726393862E-05"

$ws.Range("G9").Value = [double]"0.01352566666666667"
$ws.Range("H9").Value = [double]"0.040577"
$ws.Range("I9").Value = [double]"0.0001954433790098242"
$ws.Range("J9").Value = [double]"0.0001954433790098242"
$ws.Range("M9").Value = [double]"2.527995333333333"
$ws.Range("N9").Value = [double]"7.583985999999999"
$ws.Range("O9").Value = [double]"0.5460625962847514"
$ws.Range("P9").Value = [double]"0.5460625962847514"
$ws.Range("Q9").Value = [double]"0.03419282221355555"
$ws.Range("R9").Value = [double]"0.307735399922"
$ws.Range("S9").Value = [double]"0.0001067243189687693"
$ws.Range("T9").Value = [double]"0.0001067243189687693"

$ws.Range("G10").Value = [double]"69.18337766666667"
$ws.Range("H10").Value = [double]"207.550133"
$ws.Range("I10").Value = [double]"0.9996869977440035"
$ws.Range("J10").Value = [double]"0.9996869977440034"
$ws.Range("M10").Value = [double]"0.1293233333333333"
$ws.Range("N10").Value = [double]"0.38797"
$ws.Range("O10").Value = [double]"0.02793463826022293"
$ws.Range("P10").Value = [double]"0.02793463826022293"
$ws.Range("Q10").Value = [double]"8.947025011112224"
$ws.Range("R10").Value = [double]"80.52322510001001"
$ws.Range("S10").Value = [double]"0.02792589465542703"
$ws.Range("T10").Value = [double]"0.02792589465542703"

$ws.Range("G11").Value = [double]"69.18337766666667"
$ws.Range("H11").Value = [double]"207.550133"
$ws.Range("I11").Value = [double]"0.9996869977440035"
$ws.Range("J11").Value = [double]"0.9996869977440034"
$ws.Range("M11").Value = [double]"0.9720173333333334"
$ws.Range("N11").Value = [double]"2.916052"
$ws.Range("O11").Value = [double]"0.2099617438667928"
$ws.Range("P11").Value = [double]"0.2099617438667928"
$ws.Range("Q11").Value = [double]"67.24744227054623"
$ws.Range("R11").Value = [double]"605.2269804349161"
$ws.Range("S11").Value = [double]"0.2098960253672895"
$ws.Range("T11").Value = [double]"0.2098960253672895"

$ws.Range("G12").Value = [double]"69.18337766666667"
$ws.Range("H12").Value = [double]"207.550133"
$ws.Range("I12").Value = [double]"0.9996869977440035"
$ws.Range("J12").Value = [double]"0.9996869977440034"
$ws.Range("M12").Value = [double]"1.000161333333333"
$ws.Range("N12").Value = [double]"3.000484"
$ws.Range("O12").Value = [double]"0.2160410215882329"
$ws.Range("P12").Value = [double]"0.2160410215882329"
$ws.Range("Q12").Value = [double]"69.19453925159689"
$ws.Range("R12").Value = [double]"622.750853264372"
$ws.Range("S12").Value = [double]"0.215973400261088"
$ws.Range("T12").Value = [double]"0.215973400261088"

$ws.Range("G13").Value = [double]"69.18337766666667"
$ws.Range("H13").Value = [double]"207.550133"
$ws.Range("I13").Value = [double]"0.9996869977440035"
$ws.Range("J13").Value = [double]"0.9996869977440034"
$ws.Range("M13").Value = [double]"2.527995333333333"
$ws.Range("N13").Value = [double]"7.583985999999999"
$ws.Range("O13").Value = [double]"0.5460625962847514"
$ws.Range("P13").Value = [double]"0.5460625962847514"
$ws.Range("Q13").Value = [double]"174.8952558855709"
$ws.Range("R13").Value = [double]"1574.057302970138"
$ws.Range("S13").Value = [double]"0.545891677460199"
$ws.Range("T13").Value = [double]"0.545891677460199"
